$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# "View Search Strategy" functionality: two new "radio button" helper
# columns are inserted next to the existing Population (B) and
# Study_Types (C, soon to be D) columns, and several template /
# upload-file values are refreshed.
# ------------------------------------------------------------------

# Insert a new column before the old column C (Study_Types). After this,
# old C..G become D..H.
$ws.Columns("C:C").Insert()

# Insert a new column before the (now shifted) old column E
# (Files_to_upload, originally column E, now F after the first insert's
# shift -- the gap we want is still addressed as "E" since nothing to
# its left changed again). Old D..G (now D,F,G,H) become D,F(unchanged
# position only for D),... i.e. this places the 2nd new column right
# before the old "Files_to_upload" column.
$ws.Columns("E:E").Insert()

# --- New column C: radio-button companion to column B (Population) ---
$ws.Range("C1").Value = "Population_Radio_button"
$ws.Range("C2").Value = "Test_Sachin_2022 - Test_Sachin_radio_button"

# --- New column E: radio-button companion to column D (Study_Types) ---
$ws.Range("E1").Value = "slrtype_Radio_button"
$ws.Range("E2").Value = "Clinical_radio_button"
$ws.Range("E3").Value = "Economic_radio_button"
$ws.Range("E4").Value = "Quality of Life_radio_button"
$ws.Range("E5").Value = "Real-world Evidence_radio_button"

# Header cells in the two new columns use the default (left-aligned)
# style, unlike the rest of row 1 which is centered.
$ws.Range("C1").Style = "Normal"
$ws.Range("E1").Style = "Normal"

# --- Refresh existing data: template names and upload file paths ---
$ws.Range("F2").Value = "Clinical_search-strategy-template_Oncology.xlsx"
$ws.Range("F3").Value = "Economic_search-strategy-template_Oncology.xlsx"
$ws.Range("F4").Value = "Quality of life_search-strategy-template_Oncology.xlsx"
$ws.Range("F5").Value = "Real-world Evidence_search-strategy-template_Oncology.xlsx"

$ws.Range("G2").Value = "\Testdata\Templates\SearchStrategy\Clinical&RWE.xlsx"
$ws.Range("G3").Value = "\Testdata\Templates\SearchStrategy\ECON.xlsx"
$ws.Range("G4").Value = "\Testdata\Templates\ManageQAData\1stUpload\Cochrane RoB2 RRMM JA - Copy (3).xlsx"
$ws.Range("G5").Value = "\Testdata\Templates\ManageQAData\1stUpload\Cochrane RoB2 RRMM JA - Copy (4).xlsx"

# --- Column widths (characters); offsets chosen so the exported pixel
# width lands as close as possible to the authored widths. ---
$ws.Columns("B:B").ColumnWidth = 26.498697916666668
$ws.Columns("C:C").ColumnWidth = 38.498697916666664
$ws.Columns("E:E").ColumnWidth = 16.944010416666668

# --- View state: scrolled so column C is the leftmost visible column,
# with G3 as the active selection. ---
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("G3").Select() | Out-Null
